$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H16").Value = 3099.875
$ws.Range("I16").Value = 2949.75
$ws.Range("J16").Value = 3250
$ws.Range("K16").Value = 2949.75
$ws.Range("L16").Value = 3250
$ws.Range("M16").Value = -2719.75
$ws.Range("N16").Value = -3710
$ws.Range("H17").Value = 404.75
$ws.Range("J17").Value = 404.75
$ws.Range("L17").Value = 1214.25
$ws.Range("N17").Value = -1550.25
$ws.Range("H19").Value = 4334.3335
$ws.Range("I19").Value = 4001
$ws.Range("K19").Value = 4001
$ws.Range("H33").Value = 1316.32
$ws.Range("I33").Value = 1122.1428
$ws.Range("J33").Value = 2335.75
$ws.Range("K33").Value = 1122.1428
$ws.Range("L33").Value = 2335.75
$ws.Range("M33").Value = -893.1428000000001
$ws.Range("N33").Value = -2793.75
$ws.Range("H38").Value = 3101.5715
$ws.Range("I38").Value = 618.5
$ws.Range("J38").Value = 18000
$ws.Range("K38").Value = 1855.5
$ws.Range("L38").Value = 54000
$ws.Range("M38").Value = -1483.5
$ws.Range("N38").Value = -54744
$ws.Range("H70").Value = 4954.3335
$ws.Range("I70").Value = 4994.5
$ws.Range("K70").Value = 14983.5
$ws.Range("M70").Value = -14713.5
$ws.Range("H73").Value = 4954.3335
$ws.Range("I73").Value = 4994.5
$ws.Range("K73").Value = 14983.5
$ws.Range("M73").Value = -14047.5
$ws.Range("H98").Value = 52641556
$ws.Range("I98").Value = 62508100
$ws.Range("K98").Value = 62508100
$ws.Range("M98").Value = -62506602
$ws.Range("H122").Value = 52641556
$ws.Range("I122").Value = 62508100
$ws.Range("K122").Value = 187524300
$ws.Range("M122").Value = -187521850
$ws.Range("H138").Value = 2901.239
$ws.Range("J138").Value = 3128.805
$ws.Range("L138").Value = 9386.414999999999
$ws.Range("N138").Value = -19666.415
$ws.Range("M19").Value = -3826

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14712919
$ws.Range("I32").Value = 14712919
$ws.Range("K32").Value = 14712919
$ws.Range("M32").Value = -14712632
$ws.Range("H43").Value = 30665.5
$ws.Range("J43").Value = 33333
$ws.Range("L43").Value = 33333
$ws.Range("N43").Value = -33959
$ws.Range("H45").Value = 1910.6666
$ws.Range("I45").Value = 2091
$ws.Range("J45").Value = 1550
$ws.Range("K45").Value = 2091
$ws.Range("L45").Value = 1550
$ws.Range("M45").Value = -1714
$ws.Range("N45").Value = -2304
$ws.Range("H110").Value = 2202.5
$ws.Range("I110").Value = 1603.3334
$ws.Range("J110").Value = 4000
$ws.Range("K110").Value = 1603.3334
$ws.Range("L110").Value = 4000
$ws.Range("M110").Value = 441.6666
$ws.Range("N110").Value = -8090

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 924.71875
$ws.Range("I94").Value = 949.7
$ws.Range("K94").Value = 949.7
$ws.Range("M94").Value = -498.7

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H51").Value = 46599.9
$ws.Range("J51").Value = 70000
$ws.Range("L51").Value = 70000
$ws.Range("N51").Value = -71472
$ws.Range("H61").Value = 46599.9
$ws.Range("J61").Value = 70000
$ws.Range("L61").Value = 70000
$ws.Range("N61").Value = -70696

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 22347.389
$ws.Range("I2").Value = 111.75
$ws.Range("J2").Value = 40135.9
$ws.Range("K2").Value = 670.5
$ws.Range("L2").Value = 240815.4
$ws.Range("M2").Value = -557.5
$ws.Range("N2").Value = -241041.4
$ws.Range("H7").Value = 112.666664
$ws.Range("I7").Value = 112.666664
$ws.Range("K7").Value = 337.999992
$ws.Range("M7").Value = -225.999992
$ws.Range("H18").Value = 525
$ws.Range("I18").Value = 50
$ws.Range("J18").Value = 1000
$ws.Range("K18").Value = 150
$ws.Range("L18").Value = 3000
$ws.Range("M18").Value = 19
$ws.Range("H23").Value = 999.5
$ws.Range("I23").Value = 999
$ws.Range("K23").Value = 2997
$ws.Range("M23").Value = -2762
$ws.Range("H34").Value = 421.33334
$ws.Range("I34").Value = 421.33334
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 1264.00002
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -1180.00002
$ws.Range("H37").Value = 79993.336
$ws.Range("J37").Value = 79993.336
$ws.Range("L37").Value = 239980.008
$ws.Range("N37").Value = -240204.008
$ws.Range("H44").Value = 21765.334
$ws.Range("I44").Value = 21765.334
$ws.Range("K44").Value = 65296.00199999999
$ws.Range("M44").Value = -64898.00199999999
$ws.Range("H51").Value = 9661.091
$ws.Range("I51").Value = 5596.3335
$ws.Range("K51").Value = 16789.0005
$ws.Range("M51").Value = -16329.0005
$ws.Range("H80").Value = 4048.9092
$ws.Range("I80").Value = 2995.5
$ws.Range("K80").Value = 8986.5
$ws.Range("M80").Value = -8050.5
$ws.Range("H82").Value = 6657
$ws.Range("I82").Value = 5821.25
$ws.Range("K82").Value = 17463.75
$ws.Range("M82").Value = -17057.75
$ws.Range("H83").Value = 4048.9092
$ws.Range("I83").Value = 2995.5
$ws.Range("K83").Value = 26959.5
$ws.Range("M83").Value = -22279.5
$ws.Range("H85").Value = 6657
$ws.Range("I85").Value = 5821.25
$ws.Range("K85").Value = 17463.75
$ws.Range("M85").Value = -16059.75
$ws.Range("H97").Value = 1339.6364
$ws.Range("I97").Value = 925
$ws.Range("J97").Value = 1495.125
$ws.Range("K97").Value = 2775
$ws.Range("L97").Value = 4485.375
$ws.Range("M97").Value = -2279
$ws.Range("N97").Value = -5477.375
$ws.Range("H103").Value = 3652.5881
$ws.Range("J103").Value = 4340.857
$ws.Range("L103").Value = 13022.571
$ws.Range("N103").Value = -14780.571
$ws.Range("H113").Value = 1355.4706
$ws.Range("I113").Value = 543.75
$ws.Range("J113").Value = 1605.2307
$ws.Range("K113").Value = 1631.25
$ws.Range("L113").Value = 4815.6921
$ws.Range("M113").Value = 538.75
$ws.Range("N113").Value = -9155.6921
$ws.Range("H129").Value = 1048.2858
$ws.Range("I129").Value = 556.3333
$ws.Range("J129").Value = 4000
$ws.Range("K129").Value = 1668.9999
$ws.Range("L129").Value = 12000
$ws.Range("M129").Value = 3331.0001
$ws.Range("N129").Value = -22000
$ws.Range("N18").Value = -3338
$ws.Range("N34").ClearContents()

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2351.889
$ws.Range("I102").Value = 2400.08
$ws.Range("J102").Value = 1749.5
$ws.Range("K102").Value = 2400.08
$ws.Range("L102").Value = 1749.5
$ws.Range("M102").Value = -778.0799999999999
$ws.Range("N102").Value = -4993.5
$ws.Range("H122").Value = 1923.9286
$ws.Range("I122").Value = 1858
$ws.Range("K122").Value = 5574
$ws.Range("M122").Value = -3124
$ws.Range("H126").Value = 1934.5454
$ws.Range("I126").Value = 1989
$ws.Range("K126").Value = 5967
$ws.Range("M126").Value = -3497
$ws.Range("H132").Value = 166670600
$ws.Range("I132").Value = 200003700
$ws.Range("K132").Value = 600011100
$ws.Range("M132").Value = -600008570

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 15466499
$ws.Range("H22").Value = 3542
$ws.Range("J22").Value = 1500
$ws.Range("L22").Value = 1500
$ws.Range("N22").Value = -2090
$ws.Range("H27").Value = 3542
$ws.Range("J27").Value = 1500
$ws.Range("L27").Value = 1500
$ws.Range("N27").Value = -1714
$ws.Range("H40").Value = 2788.8235
$ws.Range("I40").Value = 1800.9231
$ws.Range("J40").Value = 5999.5
$ws.Range("K40").Value = 1800.9231
$ws.Range("L40").Value = 5999.5
$ws.Range("M40").Value = -1664.9231
$ws.Range("N40").Value = -6271.5
$ws.Range("H82").Value = 2005.5
$ws.Range("J82").Value = 2000
$ws.Range("L82").Value = 2000
$ws.Range("N82").Value = -2722
$ws.Range("H85").Value = 2005.5
$ws.Range("J85").Value = 2000
$ws.Range("L85").Value = 2000
$ws.Range("N85").Value = -4496
$ws.Range("H122").Value = 4117.8203
$ws.Range("I122").Value = 3527.9688
$ws.Range("K122").Value = 10583.9064
$ws.Range("M122").Value = -8133.9064
$ws.Range("H126").Value = 15466499

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H52").Value = 36274.2
$ws.Range("I52").Value = 35343
$ws.Range("J52").Value = 39999
$ws.Range("K52").Value = 35343
$ws.Range("L52").Value = 39999
$ws.Range("M52").Value = -35117
$ws.Range("N52").Value = -40451
$ws.Range("H107").Value = 1387.75
$ws.Range("I107").Value = 1430.9546
$ws.Range("K107").Value = 4292.8638
$ws.Range("M107").Value = -2372.8638
$ws.Range("H122").Value = 7498.3125
$ws.Range("I122").Value = 6114.8335
$ws.Range("J122").Value = 8328.4
$ws.Range("K122").Value = 18344.5005
$ws.Range("L122").Value = 24985.2
$ws.Range("M122").Value = -15894.5005
$ws.Range("N122").Value = -29885.2
